# Update salt prices and loadings in all scenarios
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 8: Magnesium chloride unit price ---
# Baseline price updated; Lower/Upper bounds are now typed-in literal values
# (no longer computed from the baseline via formula).
$ws.Range("E8").Value = 0.38
$ws.Range("G8").Value = 0.349
$ws.Range("I8").Value = 0.411

# --- Row 9: Zinc sulfate unit price ---
$ws.Range("E9").Value = 0.795
$ws.Range("G9").Value = 0.657
$ws.Range("I9").Value = 0.931

# --- Row 17: Fermentation magnesium chloride loading ---
# Lower-bound formula multiplier corrected from 0.08 to 0.8
$ws.Range("G17").Formula = "=E17*0.8"

# --- Row 18: Fermentation zinc sulfate loading ---
$ws.Range("G18").Formula = "=E18*0.8"

# Update the worksheet selection to reflect the last-edited rows (17:18)
$ws.Range("A18").Select()
